$wb = $excel.ActiveWorkbook

# --- samples_retained sheet ---
$ws1 = $wb.Worksheets.Item("samples_retained")

# Row 14 (EmoV-DB_sorted [en]) gains counts + total, and its notes are updated
$ws1.Range("C14").Value = 1317
$ws1.Range("D14").Value = 2287
$ws1.Range("E14").Value = 1568
$ws1.Range("G14").Formula = '=IF(OR(ISBLANK(C14), ISBLANK(D14),ISBLANK(E14)), "", SUM(C14:E14))'
$ws1.Range("H14").Value = "elicitation prompts based on CMU Arctic (en) and SIWIS (fr); french samples are missing :(; only 3 emos available for josh"

# --- positive sheet ---
$ws2 = $wb.Worksheets.Item("positive")

# rows 6 & 7 notes corrected from "int" to "cur"
$ws2.Range("C6").Value = "cur"
$ws2.Range("C7").Value = "cur"

# new row 10: amused / en / amu
$ws2.Range("A10").Value = "amused"
$ws2.Range("B10").Value = "en"
$ws2.Range("C10").Value = "amu"

# --- selections (set last on each sheet so the final .Select() on
#     samples_retained leaves it as the active/visible tab) ---
$ws2.Range("A11").Select()

$ws3 = $wb.Worksheets.Item("negative")
$ws3.Range("B18").Select()

$ws1.Range("C15").Select()
